$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '30.387.66'
$ws.Range("E2").Value = "'" + '  +0.80%  '

$ws.Range("D3").Value = "'" + '2.110.85'
$ws.Range("E3").Value = "'" + '  +2.53%  '

$ws.Range("D4").Value = "'" + '1.005'
$ws.Range("E4").Value = "'" + '  -0.04%  '

$ws.Range("D5").Value = "'" + '334.22'
$ws.Range("E5").Value = "'" + '  +2.88%  '

$ws.Range("D6").Value = "'" + '1.003'
$ws.Range("E6").Value = "'" + '  -0.07%  '

$ws.Range("D7").Value = "'" + '0.5234'
$ws.Range("E7").Value = "'" + '  +1.74%  '

$ws.Range("D8").Value = "'" + '0.4561'
$ws.Range("E8").Value = "'" + '  +6.35%  '

$ws.Range("D9").Value = "'" + '53.44'
$ws.Range("E9").Value = "'" + '  +17.71%  '

$ws.Range("D10").Value = "'" + '0.08920'
$ws.Range("E10").Value = "'" + '  +2.97%  '

$ws.Range("D11").Value = "'" + '1.180'
$ws.Range("E11").Value = "'" + '  +3.09%  '

$ws.Range("D12").Value = "'" + '24.44'
$ws.Range("E12").Value = "'" + '  +2.34%  '

$ws.Range("D13").Value = "'" + '2.099.84'
$ws.Range("E13").Value = "'" + '  +1.85%  '

$ws.Range("D14").Value = "'" + '6.839'
$ws.Range("E14").Value = "'" + '  +4.00%  '

$ws.Range("D15").Value = "'" + '8.048'
$ws.Range("E15").Value = "'" + '  +6.33%  '

$ws.Range("D16").Value = "'" + '96.62'
$ws.Range("E16").Value = "'" + '  +2.57%  '

$ws.Range("D17").Value = "'" + '1.004'
$ws.Range("E17").Value = "'" + '  -0.09%  '

$ws.Range("D18").Value = "'" + '0.00001136'
$ws.Range("E18").Value = "'" + '  +2.30%  '

$ws.Range("D19").Value = "'" + '0.06642'
$ws.Range("E19").Value = "'" + '  +0.74%  '

$ws.Range("D20").Value = "'" + '19.30'
$ws.Range("E20").Value = "'" + '  +3.90%  '

$ws.Range("E21").Value = "'" + '  -0.08%  '

$ws.Range("D22").Value = "'" + '6.374'

$ws.Range("D23").Value = "'" + '30.485.20'
$ws.Range("E23").Value = "'" + '  +0.85%  '

$ws.Range("D24").Value = "'" + '12.43'
$ws.Range("E24").Value = "'" + '  +2.94%  '

$ws.Range("D25").Value = "'" + '2.365'
$ws.Range("E25").Value = "'" + '  +4.15%  '

$ws.Range("D26").Value = "'" + '2.349.50'
$ws.Range("E26").Value = "'" + '  +1.95%  '

$ws.Range("E27").Value = "'" + '  +1.90%  '

$ws.Range("D28").Value = "'" + '2.572'
$ws.Range("E28").Value = "'" + '  +4.18%  '

$ws.Range("D29").Value = "'" + '163.78'
$ws.Range("E29").Value = "'" + '  +0.84%  '

$ws.Range("D30").Value = "'" + '132.86'
$ws.Range("E30").Value = "'" + '  +2.12%  '

$ws.Range("D31").Value = "'" + '1.244'
$ws.Range("E31").Value = "'" + '  +7.43%  '

$ws.Range("D32").Value = "'" + '1.711'
$ws.Range("E32").Value = "'" + '  +16.02%  '

$ws.Range("D33").Value = "'" + '0.1075'
$ws.Range("E33").Value = "'" + '  +1.70%  '

$ws.Range("D34").Value = "'" + '6.375'
$ws.Range("E34").Value = "'" + '  +5.99%  '

$ws.Range("D35").Value = "'" + '3.925'
$ws.Range("E35").Value = "'" + '  +2.29%  '

$ws.Range("D36").Value = "'" + '10.55'
$ws.Range("E36").Value = "'" + '  +11.33%  '

$ws.Range("D37").Value = "'" + '0.02589'
$ws.Range("E37").Value = "'" + '  +2.10%  '

$ws.Range("B38").Value = "'" + 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = "'" + 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").Value = "'" + '5.615'
$ws.Range("E38").Value = "'" + '  +4.68%  '

$ws.Range("B39").Value = "'" + 'Hedera'
$ws.Range("C39").Value = "'" + 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = "'" + '0.06835'
$ws.Range("E39").Value = "'" + '  +4.69%  '

$ws.Range("B40").Value = "'" + 'Algorand'
$ws.Range("C40").Value = "'" + 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").Value = "'" + '0.2309'
$ws.Range("E40").Value = "'" + '  +4.44%  '

$ws.Range("B41").Value = "'" + 'Aptos'
$ws.Range("C41").Value = "'" + 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Value = "'" + '12.77'
$ws.Range("E41").Value = "'" + '  +3.35%  '

$ws.Range("D42").Value = "'" + '0.6937'
$ws.Range("E42").Value = "'" + '  +5.40%  '

$ws.Range("D43").Value = "'" + '1.247'
$ws.Range("E43").Value = "'" + '  +1.30%  '

$ws.Range("D44").Value = "'" + '2.357'
$ws.Range("E44").Value = "'" + '  +8.81%  '

$ws.Range("E45").Value = "'" + '  -0.07%  '

$ws.Range("B46").Value = "'" + 'Decentraland'
$ws.Range("C46").Value = "'" + 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").Value = "'" + '0.6400'
$ws.Range("E46").Value = "'" + '  +2.92%  '

$ws.Range("B47").Value = "'" + 'EnergySwap'
$ws.Range("C47").Value = "'" + 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = "'" + '14.03'
$ws.Range("E47").Value = "'" + '  +1.74%  '

$ws.Range("D48").Value = "'" + '3.659'
$ws.Range("E48").Value = "'" + '  +2.10%  '

$ws.Range("D49").Value = "'" + '0.00000000351'
$ws.Range("E49").Value = "'" + '  +25.28%  '

$ws.Range("D50").Value = "'" + '1.249'
$ws.Range("E50").Value = "'" + '  +2.01%  '

$ws.Range("D51").Value = "'" + '0.3436'
$ws.Range("E51").Value = "'" + '  +28.39%  '
